$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.737660765647888
$ws.Range("B1").Value = 2.021423578262329
$ws.Range("C1").Value = 2.222915649414062
$ws.Range("D1").Value = 2.355828046798706
$ws.Range("E1").Value = 2.991458415985107
